# Capstone - Review Ranker: "Updating flask app slide"
#
# Slide 15, shape "Freeform 13" currently reads:
#   "User has to copy and paste the review into the website"
#
# The phrase "the review" is changed to "the product URL " so the final
# text becomes:
#   "User has to copy and paste the product URL into the website"
#
# The author typed this as several distinct runs (visible in the OOXML
# as separate <a:r> elements), so we rebuild the text the same way -
# splitting the paragraph's single run into five runs with the same
# boundaries:
#   1) "User has to copy and paste "
#   2) "the "
#   3) "product URL"
#   4) " "
#   5) "into the website"

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(15)
$shp = $s.Shapes.Item("Freeform 13")
$tr = $shp.TextFrame.TextRange

$full = $tr.Text

# Locate "the review" inside the current sentence so the edit is robust
# even if surrounding text/whitespace shifts slightly.
$idxThe = $full.IndexOf("the review") + 1

# Run 2: "the " - carve it out of the original single run.
$runThe = $tr.Characters($idxThe, 4)
$runThe.Text = "the "

# Run 3: "review" -> "product URL".
$idxReview = $idxThe + 4
$runReview = $tr.Characters($idxReview, 6)
$runReview.Text = "product URL"

# Run 4: the single space between "product URL" and "into".
$idxSpace = $idxReview + ("product URL").Length
$runSpace = $tr.Characters($idxSpace, 1)
$runSpace.Text = " "

# Runs 1 ("User has to copy and paste ") and 5 ("into the website") keep
# their original formatting/run boundaries automatically, since they were
# never touched above.
